# [Excel] Replace "#" by "##" in "annotation"
#
# The "fwk_content" worksheet has a column named "annotation" (column F)
# whose values are short Markdown-ish notes using "# Heading" style
# headers. This script promotes every such single "#" header marker to a
# "##" header marker, line by line, for every cell in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fwk_content")

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

# Locate the "annotation" column dynamically (defensive: falls back to F/6).
$annotationCol = 6
$headerRange = $ws.Cells.Item(1, 1)
$colCount = $usedRange.Columns.Count
for ($c = 1; $c -le $colCount; $c++) {
    $headerVal = $ws.Cells.Item(1, $c).Value2
    if ($headerVal -eq "annotation") {
        $annotationCol = $c
    }
}

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $annotationCol)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $lines = $text -split "`n"
        $newLines = @()
        foreach ($line in $lines) {
            if ($line.StartsWith("# ")) {
                $newLines += "#" + $line
            } else {
                $newLines += $line
            }
        }
        $newText = $newLines -join "`n"

        if ($newText -ne $text) {
            $cell.Value = $newText
        }
    }
}
